# Update the sleep diary date column (A2:A39) so the diary now starts on
# 2020-01-01 instead of 2000-01-01, then leave the newly-filled range
# (A2:A39) selected, matching the author's workflow of typing the new
# start date and dragging/filling it down the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startSerial = 43831  # 2020-01-01 as an Excel date serial number
$rowCount = 38         # rows 2 .. 39

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $startSerial + $i
}

$ws.Range("A2:A39").Select()
